$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Sort the "Periodo Mora" column (E16:E23) into ascending order:
# 1607, 1608, 1609, 1610, 1611, 1612, 1701, 1702
$ws.Range("E16").Value = "1607"
$ws.Range("E17").Value = "1608"
$ws.Range("E18").Value = "1609"
$ws.Range("E19").Value = "1610"
$ws.Range("E20").Value = "1611"
$ws.Range("E21").Value = "1612"
$ws.Range("E22").Value = "1701"
$ws.Range("E23").Value = "1702"
